$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10: A10 changes style from the "date-only" style (s=3) to the
# "date-time" style (s=2) used by the other data rows (same as A2's style)
$ws.Range("A10").NumberFormat = $ws.Range("A9").NumberFormat

# New row 11: A11 keeps the "date-only" style (s=3) that A10 used to have,
# B11 carries the same TotalProfit value as B10
$ws.Range("A11").NumberFormat = "YYYY-MM-DD"
$ws.Range("A11").Value = 44522
$ws.Range("B11").Value = -138.0499999999997
